# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 84
$ws1.Range("F4").Value = 253
$ws1.Range("F6").Value = 10069
$ws1.Range("F7").Value = 332
$ws1.Range("F8").Value = 908
$ws1.Range("F10").Value = 5551
$ws1.Range("F11").Value = 10
$ws1.Range("F12").Value = 22
$ws1.Range("F13").Value = 182
$ws1.Range("F15").Value = 3082
$ws1.Range("F18").Value = 595
$ws1.Range("F20").Value = 16
$ws1.Range("F22").Value = 19
$ws1.Range("F23").Value = 1529

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 84
$ws4.Range("F5").Value = 254
$ws4.Range("F7").Value = 10069
$ws4.Range("F8").Value = 332
$ws4.Range("F9").Value = 908
$ws4.Range("F11").Value = 5551
$ws4.Range("F12").Value = 10
$ws4.Range("F13").Value = 22
$ws4.Range("F14").Value = 182
$ws4.Range("F16").Value = 3082
$ws4.Range("F19").Value = 595
$ws4.Range("F21").Value = 16
$ws4.Range("F23").Value = 19
$ws4.Range("F24").Value = 1529
